$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The numeric-looking columns (runs/balls/4s/6s/sr) are stored as text,
# keep them formatted as text so the copied values don't turn numeric.
$ws.Range("G2:K2").NumberFormat = "@"

# Move the "Nov 1 2020" match (currently row 3) up into row 2,
# overwriting the "Oct 30 2020" match, then drop the now-trailing rows.
$ws.Range("A2:K2").Value2 = $ws.Range("A3:K3").Value2

# Remove old rows 3,4,5 (Oct 24 2020 and Oct 4 2020 matches, plus the
# now-duplicated Nov 1 2020 row) entirely.
$ws.Range("A3:K5").EntireRow.Delete()
